$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.726.27"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "3.770.21"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.50"
$ws.Range("E5").Value = "  +3.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.54"
$ws.Range("E6").Value = "  -4.25%  "

$ws.Range("D7").Value = "3.766.55"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("E10").Value = "  +3.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("E11").Value = "  -5.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.488"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.62"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000258"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "4.396.51"
$ws.Range("E15").Value = "  +2.27%  "

$ws.Range("D16").Value = "3.771.11"
$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").Value = "69.840.85"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.72"
$ws.Range("E20").Value = "  -1.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "507.19"
$ws.Range("E21").Value = "  -1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.54"
$ws.Range("E22").Value = "  +3.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.86"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.13"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.03"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000138"
$ws.Range("E28").Value = "  +24.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -1.97%  "

$ws.Range("E31").Value = "  +4.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  -5.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.02"
$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  +5.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.332"
$ws.Range("E38").Value = "  -3.91%  "

$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.83"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.42"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "422.15"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.68"
$ws.Range("E44").Value = "  -1.88%  "

$ws.Range("D45").Value = "3.030.73"
$ws.Range("E45").Value = "  -3.73%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.22"
$ws.Range("E48").Value = "  -4.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.58"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("E51").Value = "  +0.12%  "
